# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# These two sheets mirror the same underlying data, so both receive
# identical updates.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1299
    $ws.Range("F3").Value = 1697
    $ws.Range("F4").Value = 65
    $ws.Range("F5").Value = 6237
    $ws.Range("F6").Value = 73
}
